$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 409
$ws.Range("E6").Value = 102
$ws.Range("G6").Value = 24.93887530562347
$ws.Range("H6").Value = 75.06112469437653
